$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D cells are numeric price strings that must remain stored as TEXT
# (matching the original inlineStr type), not auto-converted to numbers.
# Setting NumberFormat to "@" (Text) before assigning the value forces Excel
# to keep it as a string; resetting the Style back to "Normal" afterwards
# avoids leaving a stray custom number-format style on the cell.
$dCells = @(
    @{Ref="D2"; Value="257.26"}
    @{Ref="D3"; Value="22.78"}
    @{Ref="D4"; Value="6.166"}
    @{Ref="D5"; Value="0.06065"}
    @{Ref="D6"; Value="6.717"}
    @{Ref="D8"; Value="1.352"}
    @{Ref="D9"; Value="0.7957"}
    @{Ref="D10"; Value="0.01322"}
    @{Ref="D11"; Value="0.1579"}
    @{Ref="D12"; Value="0.08045"}
    @{Ref="D13"; Value="0.03347"}
    @{Ref="D14"; Value="0.03093"}
    @{Ref="D15"; Value="0.09305"}
    @{Ref="D16"; Value="3.914"}
    @{Ref="D17"; Value="0.001690"}
    @{Ref="D18"; Value="0.04836"}
    @{Ref="D19"; Value="0.006182"}
    @{Ref="D21"; Value="0.003381"}
    @{Ref="D23"; Value="3.688"}
    @{Ref="D24"; Value="2.264"}
    @{Ref="D26"; Value="0.1226"}
    @{Ref="D27"; Value="0.0003018"}
    @{Ref="D40"; Value="0.04567"}
    @{Ref="D41"; Value="0.007107"}
    @{Ref="D42"; Value="0.003904"}
    @{Ref="D43"; Value="0.1112"}
    @{Ref="D44"; Value="0.009922"}
    @{Ref="D45"; Value="0.002972"}
    @{Ref="D46"; Value="0.00005910"}
    @{Ref="D48"; Value="0.7507"}
    @{Ref="D49"; Value="0.06712"}
    @{Ref="D50"; Value="0.00001501"}
)
foreach ($item in $dCells) {
    $cell = $ws.Range($item.Ref)
    $cell.NumberFormat = "@"
    $cell.Value = $item.Value
    $cell.Style = "Normal"
}

# Non-numeric text cells (coin name / link / volume label columns) can be
# assigned directly; Excel keeps them as text automatically.
$textCells = @(
    @{Ref="B10"; Value="One"}
    @{Ref="C10"; Value="https://coinranking.com/coin/6Lga5NiXX3rT+one-one"}
    @{Ref="E10"; Value="9OneONE"}
    @{Ref="B11"; Value="WazirX"}
    @{Ref="C11"; Value="https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"}
    @{Ref="E11"; Value="10WazirXWRX"}
    @{Ref="B12"; Value="MandalaExchangeToken"}
    @{Ref="C12"; Value="https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"}
    @{Ref="E12"; Value="11MandalaExchangeTokenMDX"}
    @{Ref="B13"; Value="LiechtensteinCryptoassetsExchange"}
    @{Ref="C13"; Value="https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"}
    @{Ref="E13"; Value="12LiechtensteinCryptoassetsExchangeLCX"}
    @{Ref="B14"; Value="BitrueCoin"}
    @{Ref="C14"; Value="https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"}
    @{Ref="E14"; Value="13BitrueCoinBTR"}
    @{Ref="B15"; Value="BitMartToken"}
    @{Ref="C15"; Value="https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"}
    @{Ref="E15"; Value="14BitMartTokenBMX"}
    @{Ref="B16"; Value="MCDex"}
    @{Ref="C16"; Value="https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"}
    @{Ref="E16"; Value="15MCDexMCB"}
    @{Ref="B17"; Value="BitForexToken"}
    @{Ref="C17"; Value="https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"}
    @{Ref="E17"; Value="16BitForexTokenBF"}
    @{Ref="B18"; Value="CoinExToken"}
    @{Ref="C18"; Value="https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"}
    @{Ref="E18"; Value="17CoinExTokenCET"}
    @{Ref="B42"; Value="CEJI"}
    @{Ref="C42"; Value="https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"}
    @{Ref="E42"; Value="41CEJICEJI"}
    @{Ref="B43"; Value="BKEXToken"}
    @{Ref="C43"; Value="https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"}
    @{Ref="E43"; Value="42BKEXTokenBKK"}
    @{Ref="E49"; Value="48BOLOBOLOWorstin24h"}
)
foreach ($item in $textCells) {
    $ws.Range($item.Ref).Value = $item.Value
}
